$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.0292345
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.4428765120700495
$ws.Range("J2").Value = 0.346386487911515
$ws.Range("M2").Value = 0.976531
$ws.Range("N2").Value = 1.953062
$ws.Range("O2").Value = 0.3681359341666424
$ws.Range("P2").Value = 0.3330144816631303
$ws.Range("Q2").Value = 0.0285483955195
$ws.Range("R2").Value = 0.114193582078
$ws.Range("S2").Value = 0.1630387584913719
$ws.Range("T2").Value = 0.1153517167269653

# Row 3
$ws.Range("G3").Value = 0.0292345
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.4428765120700495
$ws.Range("J3").Value = 0.346386487911515
$ws.Range("O3").Value = 0.09773758465004194
$ws.Range("P3").Value = 0.1326196171188222
$ws.Range("Q3").Value = 0.007579404683833332
$ws.Range("R3").Value = 0.045476428103
$ws.Range("S3").Value = 0.04328568058796178
$ws.Range("T3").Value = 0.04593764340195865

# Row 4
$ws.Range("G4").Value = 0.0292345
$ws.Range("H4").Value = 0.058469
$ws.Range("I4").Value = 0.4428765120700495
$ws.Range("J4").Value = 0.346386487911515
$ws.Range("M4").Value = 0.127556
$ws.Range("N4").Value = 0.382668
$ws.Range("O4").Value = 0.04808648902959583
$ws.Range("P4").Value = 0.06524830531189832
$ws.Range("Q4").Value = 0.003729035882
$ws.Range("R4").Value = 0.022374215292
$ws.Range("S4").Value = 0.0212963765391221
$ws.Range("T4").Value = 0.02260113131916671

# Row 5
$ws.Range("G5").Value = 0.0292345
$ws.Range("H5").Value = 0.058469
$ws.Range("I5").Value = 0.4428765120700495
$ws.Range("J5").Value = 0.346386487911515
$ws.Range("M5").Value = 1.116584
$ws.Range("N5").Value = 2.233168
$ws.Range("O5").Value = 0.420933584203191
$ws.Range("P5").Value = 0.380775051681252
$ws.Range("Q5").Value = 0.032642774948
$ws.Range("R5").Value = 0.130571099792
$ws.Range("S5").Value = 0.1864215975850537
$ws.Range("T5").Value = 0.1318953328361945

# Row 6
$ws.Range("G6").Value = 0.0292345
$ws.Range("H6").Value = 0.058469
$ws.Range("I6").Value = 0.4428765120700495
$ws.Range("J6").Value = 0.346386487911515
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03536566666666666
$ws.Range("N6").Value = 0.106097
$ws.Range("O6").Value = 0.01333226772704545
$ws.Range("P6").Value = 0.01809048430670052
$ws.Range("Q6").Value = 0.001033897582166667
$ws.Range("R6").Value = 0.006203385493
$ws.Range("S6").Value = 0.005904548228937974
$ws.Range("T6").Value = 0.006266299323616373

# Row 7
$ws.Range("G7").Value = 0.0292345
$ws.Range("H7").Value = 0.058469
$ws.Range("I7").Value = 0.4428765120700495
$ws.Range("J7").Value = 0.346386487911515
$ws.Range("M7").Value = 0.137338
$ws.Range("N7").Value = 0.412014
$ws.Range("O7").Value = 0.05177414022348326
$ws.Range("P7").Value = 0.07025205991819664
$ws.Range("Q7").Value = 0.004015007760999999
$ws.Range("R7").Value = 0.024090046566
$ws.Range("S7").Value = 0.02292955063760192
$ws.Range("T7").Value = 0.02433436430361345

# Row 8
$ws.Range("I8").Value = 0.5571234879299505
$ws.Range("J8").Value = 0.6536135120884849
$ws.Range("M8").Value = 0.976531
$ws.Range("N8").Value = 1.953062
$ws.Range("O8").Value = 0.3681359341666424
$ws.Range("P8").Value = 0.3330144816631303
$ws.Range("Q8").Value = 0.035912904056
$ws.Range("R8").Value = 0.215477424336
$ws.Range("S8").Value = 0.2050971756752704
$ws.Range("T8").Value = 0.217662764936165

# Row 9
$ws.Range("I9").Value = 0.5571234879299505
$ws.Range("J9").Value = 0.6536135120884849
$ws.Range("O9").Value = 0.09773758465004194
$ws.Range("P9").Value = 0.1326196171188222
$ws.Range("S9").Value = 0.05445190406208016
$ws.Range("T9").Value = 0.08668197371686354

# Row 10
$ws.Range("I10").Value = 0.5571234879299505
$ws.Range("J10").Value = 0.6536135120884849
$ws.Range("M10").Value = 0.127556
$ws.Range("N10").Value = 0.382668
$ws.Range("O10").Value = 0.04808648902959583
$ws.Range("P10").Value = 0.06524830531189832
$ws.Range("Q10").Value = 0.004690999455999999
$ws.Range("R10").Value = 0.042218995104
$ws.Range("S10").Value = 0.02679011249047373
$ws.Range("T10").Value = 0.04264717399273161

# Row 11
$ws.Range("I11").Value = 0.5571234879299505
$ws.Range("J11").Value = 0.6536135120884849
$ws.Range("M11").Value = 1.116584
$ws.Range("N11").Value = 2.233168
$ws.Range("O11").Value = 0.420933584203191
$ws.Range("P11").Value = 0.380775051681252
$ws.Range("Q11").Value = 0.041063493184
$ws.Range("R11").Value = 0.246380959104
$ws.Range("S11").Value = 0.2345119866181373
$ws.Range("T11").Value = 0.2488797188450575

# Row 12
$ws.Range("I12").Value = 0.5571234879299505
$ws.Range("J12").Value = 0.6536135120884849
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.03536566666666666
$ws.Range("N12").Value = 0.106097
$ws.Range("O12").Value = 0.01333226772704545
$ws.Range("P12").Value = 0.01809048430670052
$ws.Range("Q12").Value = 0.001300607757333333
$ws.Range("R12").Value = 0.011705469816
$ws.Range("S12").Value = 0.007427719498107473
$ws.Range("T12").Value = 0.01182418498308415

# Row 13
$ws.Range("I13").Value = 0.5571234879299505
$ws.Range("J13").Value = 0.6536135120884849
$ws.Range("M13").Value = 0.137338
$ws.Range("N13").Value = 0.412014
$ws.Range("O13").Value = 0.05177414022348326
$ws.Range("P13").Value = 0.07025205991819664
$ws.Range("Q13").Value = 0.007579404683833332
$ws.Range("R13").Value = 0.045456680592
$ws.Range("S13").Value = 0.02884458958588134
$ws.Range("T13").Value = 0.04591769561458318
